# Generate Report for Handoff
# Rotates the localization-status report onto a new source file
# (57efe6cb-9d3c-4b27-94ab-05db38a15e75.md) and a fresh handoff round:
# new handoff xliff files/timestamps, and the not-yet-handed-back
# target/handback columns reset to blank / the zero DateTime sentinel.

$wb = $excel.ActiveWorkbook

$oldGuid = "d3e4c626-60ad-4f37-8ad9-8d68936cf614"
$newGuid = "57efe6cb-9d3c-4b27-94ab-05db38a15e75"

$oldHash = "ebf543cb0225f370037be2c7db637e8e5e101298"
$newHash = "db6395e9f575822395d88754a7130dca2312501d"

$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 02:53:46"

# Rebuild the B2 hyperlink so its display text follows the new file name
# while keeping the same external target (only the file name inside the
# URL changes).
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/988743f926ef7ebe5c80ff4f90b757c72d0b5625/e2e/$newGuid.md",
    "",
    "",
    "e2e\$newGuid.md"
)

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-16 02:53:41"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $zeroDate

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# Drop every hyperlink on the sheet (A2 + I2) then recreate only the A2
# one, pointed at the new file name - I2 no longer carries a hyperlink
# since "Latest Target File" is now blank.
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/988743f926ef7ebe5c80ff4f90b757c72d0b5625/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-16 02:53:46"
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $zeroDate

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/988743f926ef7ebe5c80ff4f90b757c72d0b5625/e2e/$newGuid.md",
    "",
    "",
    "$newGuid.md"
)
